# Gallery_AttachedFunctionality_FIM_Node.xlsx - add "Turkey" market test data sheet
# (copied from the "Spain" sheet, same layout, new market name + JIRA ref).

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Make sure Spain is the active sheet and its original A1:D11 block is the
# selection before we branch off a copy - the copy inherits the source
# sheet's current view/selection.
$spain.Activate()
$spain.Range("A1:D11").Select()

# Duplicate "Spain" right after itself; the copy becomes the new active sheet.
$spain.Copy($null, $spain)

$turkey = $wb.Worksheets.Item("Spain (2)")
$turkey.Name = "Turkey"

# Update the market name and JIRA/ticket reference for Turkey.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3322/T3308"

# Column widths specific to the Turkey sheet.
$turkey.Columns.Item(1).ColumnWidth = 24.33
$turkey.Columns.Item(2).ColumnWidth = 20.33
$turkey.Columns.Item(4).ColumnWidth = 20.17

# Leave Turkey as the active/selected tab with D17 selected.
$turkey.Activate()
$turkey.Range("D17").Select()
